$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: merge the two "SAT Jan 05" / " 11:35:14 IST 2019" runs into
# a single run by doing an identity Find & Replace (formatting unchanged).
# ---------------------------------------------------------------------
$null = $d.Content.Find.Execute("SAT Jan 05 11:35:14 IST 2019", $false, $false, $false, $false, $false, $true, 1, $false, "SAT Jan 05 11:35:14 IST 2019", 2)

# ---------------------------------------------------------------------
# Change 2: append a new purchase-details block (MAMATHA CHICK IN,
# 10/01/2019) after the very last "- CASH AND CLEARD" paragraph.
# ---------------------------------------------------------------------

# Locate the last paragraph containing "- CASH AND CLEARD".
$count = $d.Paragraphs.Count
$lastIdx = -1
for ($i = $count; $i -ge 1; $i--) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*- CASH AND CLEARD*") {
        $lastIdx = $i
        break
    }
}

$cur = $d.Paragraphs.Item($lastIdx)

# 1) blank paragraph
$cur.Range.InsertParagraphAfter()
$lastIdx = $lastIdx + 1
$cur = $d.Paragraphs.Item($lastIdx)

# 2) "TUE Jan 08" / " 11:15:57 IST 2019" -- kept as two distinct runs.
$cur.Range.InsertParagraphAfter()
$lastIdx = $lastIdx + 1
$cur = $d.Paragraphs.Item($lastIdx)
$cur.Range.InsertAfter("TUE Jan 08")

$cur2 = $d.Paragraphs.Item($lastIdx)
$cur2.Range.InsertParagraphAfter()
$tmpIdx = $lastIdx + 1
$tmpPara = $d.Paragraphs.Item($tmpIdx)
$tmpPara.Range.InsertAfter(" 11:15:57 IST 2019")
$endOfFirst = $d.Paragraphs.Item($lastIdx).Range.End
$sep = $d.Range($endOfFirst - 1, $endOfFirst)
$sep.Delete()

$cur = $d.Paragraphs.Item($lastIdx)

# 3) Person Name ... - CHANDU
$cur.Range.InsertParagraphAfter()
$lastIdx = $lastIdx + 1
$cur = $d.Paragraphs.Item($lastIdx)
$cur.Range.InsertAfter("Person Name`t`t`t`t- CHANDU")

# 4) Bill number ... - 10060
$cur.Range.InsertParagraphAfter()
$lastIdx = $lastIdx + 1
$cur = $d.Paragraphs.Item($lastIdx)
$cur.Range.InsertAfter("Bill number`t`t`t`t- 10060")

# 5) dashed separator line
$cur.Range.InsertParagraphAfter()
$lastIdx = $lastIdx + 1
$cur = $d.Paragraphs.Item($lastIdx)
$cur.Range.InsertAfter("---------------------------------------------------------------")

# 6) Item Name ... - SORE KAI
$cur.Range.InsertParagraphAfter()
$lastIdx = $lastIdx + 1
$cur = $d.Paragraphs.Item($lastIdx)
$cur.Range.InsertAfter("Item Name`t`t`t`t- SORE KAI")

# 7) Number of Pockets ... - 2
$cur.Range.InsertParagraphAfter()
$lastIdx = $lastIdx + 1
$cur = $d.Paragraphs.Item($lastIdx)
$cur.Range.InsertAfter("Number of Pockets`t`t`t- 2")

# 8) Number of KGs ... - 102
$cur.Range.InsertParagraphAfter()
$lastIdx = $lastIdx + 1
$cur = $d.Paragraphs.Item($lastIdx)
$cur.Range.InsertAfter("Number of KGs`t`t`t- 102")

# 9) Rate ... - 8
$cur.Range.InsertParagraphAfter()
$lastIdx = $lastIdx + 1
$cur = $d.Paragraphs.Item($lastIdx)
$cur.Range.InsertAfter("Rate`t`t`t`t`t- 8")

# 10) Total Price ... - 816.0
$cur.Range.InsertParagraphAfter()
$lastIdx = $lastIdx + 1
$cur = $d.Paragraphs.Item($lastIdx)
$cur.Range.InsertAfter("Total Price`t`t`t`t- 816.0")

# 11) Amount balance ... - 816.0 (bold)
$cur.Range.InsertParagraphAfter()
$lastIdx = $lastIdx + 1
$cur = $d.Paragraphs.Item($lastIdx)
$cur.Range.Font.Bold = $true
$cur.Range.InsertAfter("Amount balance`t`t`t- 816.0")

# 12) blank paragraph, bold paragraph mark
$cur.Range.InsertParagraphAfter()
$lastIdx = $lastIdx + 1
$cur = $d.Paragraphs.Item($lastIdx)
$cur.Range.Font.Bold = $true

# 13) blank paragraph, normal paragraph mark
$cur.Range.InsertParagraphAfter()
$lastIdx = $lastIdx + 1
$cur = $d.Paragraphs.Item($lastIdx)
$cur.Range.Font.Bold = $false

Write-Host "done"
